$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 695.8333
$ws.Range("H57").Value = 37370
$ws.Range("J57").Value = 43226.668
$ws.Range("L57").Value = 129680.004
$ws.Range("N57").Value = -130678.004
$ws.Range("H88").Value = 5601.4165
$ws.Range("I88").Value = 1127
$ws.Range("J88").Value = 14550.25
$ws.Range("K88").Value = 1127
$ws.Range("L88").Value = 14550.25
$ws.Range("M88").Value = -721
$ws.Range("N88").Value = -15362.25
$ws.Range("H91").Value = 5601.4165
$ws.Range("I91").Value = 1127
$ws.Range("J91").Value = 14550.25
$ws.Range("K91").Value = 1127
$ws.Range("L91").Value = 14550.25
$ws.Range("M91").Value = 277
$ws.Range("N91").Value = -17358.25
$ws.Range("H138").Value = 3150.55
$ws.Range("I138").Value = 1919.25
$ws.Range("K138").Value = 5757.75
$ws.Range("M138").Value = -617.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22107.672
$ws.Range("I32").Value = 22404.316
$ws.Range("J32").Value = 5199
$ws.Range("K32").Value = 22404.316
$ws.Range("L32").Value = 5199
$ws.Range("M32").Value = -22117.316
$ws.Range("N32").Value = -5773
$ws.Range("H97").Value = 1985.8536
$ws.Range("I97").Value = 1675.7084
$ws.Range("J97").Value = 2423.7058
$ws.Range("K97").Value = 1675.7084
$ws.Range("L97").Value = 2423.7058
$ws.Range("M97").Value = -1179.7084
$ws.Range("N97").Value = -3415.7058
$ws.Range("H110").Value = 1168
$ws.Range("I110").Value = 1108.091
$ws.Range("K110").Value = 1108.091
$ws.Range("M110").Value = 936.9090000000001
$ws.Range("H122").Value = 3077.35
$ws.Range("I122").Value = 2999.7144
$ws.Range("J122").Value = 3258.5
$ws.Range("K122").Value = 8999.143199999999
$ws.Range("L122").Value = 9775.5
$ws.Range("M122").Value = -6549.143199999999
$ws.Range("N122").Value = -14675.5
$ws.Range("H132").Value = 1858.0465
$ws.Range("I132").Value = 1297.3235
$ws.Range("J132").Value = 3976.3333
$ws.Range("K132").Value = 3891.9705
$ws.Range("L132").Value = 11928.9999
$ws.Range("M132").Value = -1361.9705
$ws.Range("N132").Value = -16988.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2989.8386
$ws.Range("I134").Value = 2847.9312
$ws.Range("K134").Value = 8543.793600000001
$ws.Range("M134").Value = -6008.793600000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1589.5186
$ws.Range("I16").Value = 1530.95
$ws.Range("K16").Value = 1530.95
$ws.Range("M16").Value = -1243.95
$ws.Range("H22").Value = 291.23077
$ws.Range("I22").Value = 282.33334
$ws.Range("J22").Value = 298.85715
$ws.Range("K22").Value = 282.33334
$ws.Range("L22").Value = 298.85715
$ws.Range("M22").Value = 67.66665999999998
$ws.Range("N22").Value = -998.85715
$ws.Range("H31").Value = 5559054.5
$ws.Range("I31").Value = 5885851.5
$ws.Range("K31").Value = 5885851.5
$ws.Range("M31").Value = -5885556.5
$ws.Range("H34").Value = 5559054.5
$ws.Range("I34").Value = 5885851.5
$ws.Range("K34").Value = 5885851.5
$ws.Range("M34").Value = -5885649.5
$ws.Range("H58").Value = 2313
$ws.Range("I58").Value = 2075.2632
$ws.Range("J58").Value = 3216.4
$ws.Range("K58").Value = 2075.2632
$ws.Range("L58").Value = 3216.4
$ws.Range("M58").Value = -1872.2632
$ws.Range("N58").Value = -3622.4
$ws.Range("H62").Value = 6329.875
$ws.Range("I62").Value = 2564.3333
$ws.Range("J62").Value = 11171.286
$ws.Range("K62").Value = 2564.3333
$ws.Range("L62").Value = 11171.286
$ws.Range("M62").Value = -1940.3333
$ws.Range("N62").Value = -12419.286
$ws.Range("H65").Value = 6329.875
$ws.Range("I65").Value = 2564.3333
$ws.Range("J65").Value = 11171.286
$ws.Range("K65").Value = 12821.6665
$ws.Range("L65").Value = 55856.43
$ws.Range("M65").Value = -9701.666499999999
$ws.Range("N65").Value = -62096.43
$ws.Range("H99").Value = 6496.364
$ws.Range("I99").Value = 5046.5
$ws.Range("J99").Value = 20995
$ws.Range("K99").Value = 5046.5
$ws.Range("L99").Value = 20995
$ws.Range("M99").Value = -3548.5
$ws.Range("N99").Value = -23991
$ws.Range("H113").Value = 1589.5186
$ws.Range("I113").Value = 1530.95
$ws.Range("K113").Value = 1530.95
$ws.Range("M113").Value = 639.05
$ws.Range("H122").Value = 24221.334
$ws.Range("J122").Value = 2358
$ws.Range("L122").Value = 7074
$ws.Range("N122").Value = -11974
$ws.Range("H126").Value = 6496.364
$ws.Range("I126").Value = 5046.5
$ws.Range("J126").Value = 20995
$ws.Range("K126").Value = 15139.5
$ws.Range("L126").Value = 62985
$ws.Range("M126").Value = -12669.5
$ws.Range("N126").Value = -67925
$ws.Range("H132").Value = 31470.75
$ws.Range("I132").Value = 38575.125
$ws.Range("J132").Value = 3053.25
$ws.Range("K132").Value = 115725.375
$ws.Range("L132").Value = 9159.75
$ws.Range("M132").Value = -113195.375
$ws.Range("N132").Value = -14219.75
$ws.Range("H134").Value = 3673.8462
$ws.Range("I134").Value = 2274.1428
$ws.Range("K134").Value = 6822.428400000001
$ws.Range("M134").Value = -4287.428400000001
$ws.Range("H135").Value = 111876.25
$ws.Range("J135").Value = 111876.25
$ws.Range("L135").Value = 111876.25
$ws.Range("N135").Value = -122016.25
$ws.Range("H136").Value = 2313
$ws.Range("I136").Value = 2075.2632
$ws.Range("J136").Value = 3216.4
$ws.Range("K136").Value = 6225.7896
$ws.Range("L136").Value = 9649.200000000001
$ws.Range("M136").Value = -3675.7896
$ws.Range("N136").Value = -14749.2
$ws.Range("H138").Value = 113797.375
$ws.Range("J138").Value = 113797.375
$ws.Range("L138").Value = 113797.375
$ws.Range("N138").Value = -124077.375
$ws.Range("H139").Value = 69495
$ws.Range("J139").Value = 69495
$ws.Range("L139").Value = 69495
$ws.Range("N139").Value = -79775
$ws.Range("H140").Value = 113463.5
$ws.Range("J140").Value = 113463.5
$ws.Range("L140").Value = 113463.5
$ws.Range("N140").Value = -123823.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 8112.3335
$ws.Range("I87").Value = 2835.1667
$ws.Range("K87").Value = 8505.500100000001
$ws.Range("M87").Value = -7257.500100000001
$ws.Range("H90").Value = 8112.3335
$ws.Range("I90").Value = 2835.1667
$ws.Range("K90").Value = 25516.5003
$ws.Range("M90").Value = -19276.5003
$ws.Range("H107").Value = 1388.579
$ws.Range("I107").Value = 1725.5555
$ws.Range("K107").Value = 5176.666499999999
$ws.Range("M107").Value = -3256.666499999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 500005000
$ws.Range("I11").Value = 1000000000
$ws.Range("J11").Value = 10002
$ws.Range("K11").Value = 1000000000
$ws.Range("L11").Value = 10002
$ws.Range("M11").Value = -999999861
$ws.Range("N11").Value = -10280
$ws.Range("H18").Value = 22000
$ws.Range("J18").Value = 22000
$ws.Range("L18").Value = 22000
$ws.Range("N18").Value = -22586
$ws.Range("H107").Value = 323.1111
$ws.Range("I107").Value = 188.25
$ws.Range("J107").Value = 431
$ws.Range("K107").Value = 188.25
$ws.Range("L107").Value = 431
$ws.Range("M107").Value = 1731.75
$ws.Range("N107").Value = -4271
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7363.0454
$ws.Range("I122").Value = 7421.968
$ws.Range("J122").Value = 7222.5386
$ws.Range("K122").Value = 22265.904
$ws.Range("L122").Value = 21667.6158
$ws.Range("M122").Value = -19815.904
$ws.Range("N122").Value = -26567.6158
$ws.Range("H136").Value = 9249.5
$ws.Range("J136").Value = 9666
$ws.Range("L136").Value = 28998
$ws.Range("N136").Value = -34098
$ws.Range("H137").Value = 98990
$ws.Range("J137").Value = 98990
$ws.Range("L137").Value = 98990
$ws.Range("N137").Value = -109190

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3367.3076
$ws.Range("I17").Value = 3502.0833
$ws.Range("K17").Value = 3502.0833
$ws.Range("M17").Value = -3330.0833
$ws.Range("H74").Value = 15668.875
$ws.Range("I74").Value = 7784.5
$ws.Range("J74").Value = 18297
$ws.Range("K74").Value = 7784.5
$ws.Range("L74").Value = 18297
$ws.Range("M74").Value = -6848.5
$ws.Range("N74").Value = -20169
$ws.Range("H77").Value = 15668.875
$ws.Range("I77").Value = 7784.5
$ws.Range("J77").Value = 18297
$ws.Range("K77").Value = 23353.5
$ws.Range("L77").Value = 54891
$ws.Range("M77").Value = -18673.5
$ws.Range("N77").Value = -64251
$ws.Range("H94").Value = 40270
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 881.25
$ws.Range("I113").Value = 604
$ws.Range("J113").Value = 1343.3334
$ws.Range("K113").Value = 1812
$ws.Range("L113").Value = 4030.0002
$ws.Range("M113").Value = 358
$ws.Range("N113").Value = -8370.0002
$ws.Range("H122").Value = 291975.88
$ws.Range("I122").Value = 461861.4
$ws.Range("J122").Value = 8833.333000000001
$ws.Range("K122").Value = 1385584.2
$ws.Range("L122").Value = 26499.999
$ws.Range("M122").Value = -1383134.2
$ws.Range("N122").Value = -31399.999
$ws.Range("H136").Value = 18507.348
$ws.Range("I136").Value = 21284.375
$ws.Range("J136").Value = 6165
$ws.Range("K136").Value = 63853.125
$ws.Range("L136").Value = 18495
$ws.Range("M136").Value = -61303.125
$ws.Range("N136").Value = -23595
$ws.Range("H139").Value = 122306.5
$ws.Range("J139").Value = 122306.5
$ws.Range("L139").Value = 122306.5
$ws.Range("N139").Value = -132586.5

Write-Host "Applied all updates"